{"js": "// Replace \"Spring\" with \"Django\" in the \"Backend server: ...\" bullet of the\n// Technologies Stack list (the other two \"Spring\" occurrences in the doc -\n// the prose paragraph and the \"Spring REST\" bullet - must stay untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Backend server:\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Backend server:' paragraph\");\n}\n\n// Scope the search to just this paragraph so only this \"Spring\" is touched.\nconst hits = target.search(\"Spring\", { matchCase: true, matchWholeWord: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'Spring' match in the paragraph, found \" + hits.items.length);\n}\n\n// \"Replace\" keeps the surrounding \"Backend server: ... Flask or \" / \". \"\n// text intact and only swaps the matched word itself.\nhits.items[0].insertText(\"Django\", \"Replace\");\nawait context.sync();\n", "ps1": "# Replace \"Spring\" with \"Django\" in the \"Backend server: ...\" bullet of the\n# Technologies Stack list (the other two \"Spring\" occurrences in the doc -\n# the prose paragraph and the \"Spring REST\" bullet - must stay untouched).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Backend server:\")) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the 'Backend server:' paragraph\"\n}\n\n# Scope the Find/Replace to just this paragraph's range so only this\n# \"Spring\" is touched (the word also appears in two other paragraphs).\n$rng = $target.Range\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format, ReplaceWith,\n# Replace(1=wdReplaceOne)\n$result = $find.Execute(\"Spring\", $true, $true, $false, $false, $false, $true, 1, $false, \"Django\", 1)\nif (-not $result) {\n    throw \"Could not find 'Spring' in the 'Backend server:' paragraph\"\n}\n"}
